$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 1
$ws.Range("B7").Value = 1

$ws.Range("B7").Select()
